$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 2181.6667
$ws.Range("I40").Value = 1845
$ws.Range("J40").Value = 2212.2727
$ws.Range("K40").Value = 1845
$ws.Range("L40").Value = 2212.2727
$ws.Range("M40").Value = -1670
$ws.Range("N40").Value = -2562.2727
# row 42 (Leve Item ID 4600)
$ws.Range("H42").Value = 328.57144
$ws.Range("J42").Value = 366.66666
$ws.Range("L42").Value = 1099.99998
$ws.Range("N42").Value = -1559.99998
# row 64 (Leve Item ID 5506)
$ws.Range("H64").Value = 5282.6113
$ws.Range("I64").Value = 6935.25
$ws.Range("J64").Value = 3960.5
$ws.Range("K64").Value = 6935.25
$ws.Range("L64").Value = 3960.5
$ws.Range("M64").Value = -6687.25
$ws.Range("N64").Value = -4456.5
# row 67 (Leve Item ID 5506)
$ws.Range("H67").Value = 5282.6113
$ws.Range("I67").Value = 6935.25
$ws.Range("J67").Value = 3960.5
$ws.Range("K67").Value = 6935.25
$ws.Range("L67").Value = 3960.5
$ws.Range("M67").Value = -6077.25
$ws.Range("N67").Value = -5676.5
# row 74 (Leve Item ID 5507)
$ws.Range("H74").Value = 2998816.5
$ws.Range("I74").Value = 3185930
$ws.Range("K74").Value = 3185930
$ws.Range("M74").Value = -3184994
# row 77 (Leve Item ID 5507)
$ws.Range("H77").Value = 2998816.5
$ws.Range("I77").Value = 3185930
$ws.Range("K77").Value = 15929650
$ws.Range("M77").Value = -15924970
# row 106 (Leve Item ID 19903)
$ws.Range("H106").Value = 2589.353
$ws.Range("I106").Value = 2655.3076
$ws.Range("K106").Value = 2655.3076
$ws.Range("M106").Value = -2024.3076
# row 129 (Leve Item ID 36115)
$ws.Range("H129").Value = 1017
$ws.Range("I129").Value = 441.44446
$ws.Range("J129").Value = 1304.7778
$ws.Range("K129").Value = 1324.33338
$ws.Range("L129").Value = 3914.3334
$ws.Range("M129").Value = 3675.66662
$ws.Range("N129").Value = -13914.3334

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 64 (Leve Item ID 10664)
$ws.Range("H64").Value = 29000
$ws.Range("J64").Value = 29000
$ws.Range("L64").Value = 29000
$ws.Range("N64").Value = -29496
# row 67 (Leve Item ID 10664)
$ws.Range("H67").Value = 29000
$ws.Range("J67").Value = 29000
$ws.Range("L67").Value = 29000
$ws.Range("N67").Value = -30716
# row 76 (Leve Item ID 10679)
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
# row 79 (Leve Item ID 10679)
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
# row 88 (Leve Item ID 12530)
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()
# row 91 (Leve Item ID 12530)
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()
# row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 3556.257
$ws.Range("I102").Value = 2978.76
$ws.Range("K102").Value = 2978.76
$ws.Range("M102").Value = -1356.76
# row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 1314.159
$ws.Range("I132").Value = 850.63336
$ws.Range("K132").Value = 2551.90008
$ws.Range("M132").Value = -21.90008000000034

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 2316.5
$ws.Range("I86").Value = 2453.2
$ws.Range("J86").Value = 1633
$ws.Range("K86").Value = 2453.2
$ws.Range("L86").Value = 1633
$ws.Range("M86").Value = -1330.2
$ws.Range("N86").Value = -3879
# row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 2316.5
$ws.Range("I89").Value = 2453.2
$ws.Range("J89").Value = 1633
$ws.Range("K89").Value = 12266
$ws.Range("L89").Value = 8165
$ws.Range("M89").Value = -6650
$ws.Range("N89").Value = -19397
# row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 65424.188
$ws.Range("I134").Value = 3098.3333
$ws.Range("J134").Value = 252401.75
$ws.Range("K134").Value = 9294.999899999999
$ws.Range("L134").Value = 757205.25
$ws.Range("M134").Value = -6759.999899999999
$ws.Range("N134").Value = -762275.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2229.6428
$ws.Range("I31").Value = 2232.3948
$ws.Range("J31").Value = 2203.5
$ws.Range("K31").Value = 2232.3948
$ws.Range("L31").Value = 2203.5
$ws.Range("M31").Value = -1937.3948
$ws.Range("N31").Value = -2793.5
# row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2229.6428
$ws.Range("I34").Value = 2232.3948
$ws.Range("J34").Value = 2203.5
$ws.Range("K34").Value = 2232.3948
$ws.Range("L34").Value = 2203.5
$ws.Range("M34").Value = -2030.3948
$ws.Range("N34").Value = -2607.5
# row 81 (Leve Item ID 10613)
$ws.Range("H81").Value = 35000
$ws.Range("J81").Value = 35000
$ws.Range("L81").Value = 35000
$ws.Range("N81").Value = -36996
# row 84 (Leve Item ID 10613)
$ws.Range("H84").Value = 35000
$ws.Range("J84").Value = 35000
$ws.Range("L84").Value = 105000
$ws.Range("N84").Value = -114984

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 80 (Leve Item ID 12890)
$ws.Range("H80").Value = 7444.4443
$ws.Range("I80").Value = 3500
$ws.Range("J80").Value = 7937.5
$ws.Range("K80").Value = 10500
$ws.Range("L80").Value = 23812.5
$ws.Range("M80").Value = -9564
$ws.Range("N80").Value = -25684.5
# row 83 (Leve Item ID 12890)
$ws.Range("H83").Value = 7444.4443
$ws.Range("I83").Value = 3500
$ws.Range("J83").Value = 7937.5
$ws.Range("K83").Value = 31500
$ws.Range("L83").Value = 71437.5
$ws.Range("M83").Value = -26820
$ws.Range("N83").Value = -80797.5
# row 94 (Leve Item ID 19811)
$ws.Range("H94").Value = 4513.5
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 4513.5
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 13540.5
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -14892.5
# row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 487065.88
$ws.Range("I107").Value = 992.9
$ws.Range("J107").Value = 1297187.5
$ws.Range("K107").Value = 2978.7
$ws.Range("L107").Value = 3891562.5
$ws.Range("M107").Value = -1058.7
$ws.Range("N107").Value = -3895402.5
# row 112 (Leve Item ID 27855)
$ws.Range("H112").Value = 3845.7
$ws.Range("I112").Value = 2509
$ws.Range("J112").Value = 4418.5713
$ws.Range("K112").Value = 7527
$ws.Range("L112").Value = 13255.7139
$ws.Range("M112").Value = -6419
$ws.Range("N112").Value = -15471.7139
# row 121 (Leve Item ID 27878)
$ws.Range("H121").Value = 866
$ws.Range("I121").Value = 501.16666
$ws.Range("J121").Value = 1109.2222
$ws.Range("K121").Value = 1503.49998
$ws.Range("L121").Value = 3327.6666
$ws.Range("M121").Value = -193.4999800000001
$ws.Range("N121").Value = -5947.6666

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 4137.5
$ws.Range("I126").Value = 4211.6
$ws.Range("J126").Value = 4014
$ws.Range("K126").Value = 12634.8
$ws.Range("L126").Value = 12042
$ws.Range("M126").Value = -10164.8
$ws.Range("N126").Value = -16982

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 1149.5625
$ws.Range("I93").Value = 981.1818
$ws.Range("K93").Value = 981.1818
$ws.Range("M93").Value = 266.8182

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 86 (Leve Item ID 11977)
$ws.Range("H86").Value = 17000
$ws.Range("J86").Value = 17000
$ws.Range("L86").Value = 17000
$ws.Range("N86").Value = -19246
# row 89 (Leve Item ID 11977)
$ws.Range("H89").Value = 17000
$ws.Range("J89").Value = 17000
$ws.Range("L89").Value = 85000
$ws.Range("N89").Value = -96232
# row 93 (Leve Item ID 19613)
$ws.Range("H93").Value = 27956.357
$ws.Range("J93").Value = 27956.357
$ws.Range("L93").Value = 27956.357
$ws.Range("N93").Value = -32948.357
# row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1608.1818
$ws.Range("I132").Value = 1608.1818
$ws.Range("K132").Value = 4824.5454
$ws.Range("M132").Value = -2294.5454
# row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 1772.8276
$ws.Range("I136").Value = 1771.3334
$ws.Range("K136").Value = 5314.0002
$ws.Range("M136").Value = -2764.0002
